$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append a new sentence (as its own run, same bold/italic formatting) to
#    the paragraph that ends with "... with gene expression information."
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("gene expression information.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'gene expression information.' sentence"
}
$targetPara = $rng.Paragraphs(1)

$sentence = " The first 5 are " + [char]0x201C + "sanity check" + [char]0x201D + " variables (because we know what the pattern should be with PCs). "

# Create a fresh paragraph right after the target paragraph, fill it with the
# new sentence text (it inherits the bold/italic paragraph mark formatting),
# then delete the paragraph break that separates the two paragraphs so the
# new sentence becomes a second run inside the original paragraph.
$targetPara.Range.InsertParagraphAfter()
$sentencePara = $targetPara.Next()
$sentencePara.Range.Text = $sentence

$breakRange = $d.Range($targetPara.Range.End - 1, $targetPara.Range.End)
$breakRange.Delete()

# ---------------------------------------------------------------------------
# 2) Insert five new "ListParagraph" (numId 6) bullet items -- Sample ID,
#    Sample name, Species, Tissue, Individual -- right before the existing
#    "RNA extraction date" bullet (which is two paragraphs after the
#    paragraph updated above: the blank bold/italic paragraph, then the
#    bullet list).
# ---------------------------------------------------------------------------
$blankPara = $targetPara.Next()
$rnaPara = $blankPara.Next()

$items = @("Sample ID", "Sample name", "Species", "Tissue", "Individual")
foreach ($item in $items) {
    $rnaPara.Range.InsertParagraphBefore()
    $rnaPara.Range.Text = $item
    $rnaPara = $rnaPara.Next()
}

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark: remove it from the end of the document
#    (end of "Sequence encoding") and place it at the end of the newly
#    added "Individual" item.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$individualPara = $rnaPara.Previous()
$individualEnd = $individualPara.Range.End - 1
$bookmarkRange = $d.Range($individualEnd, $individualEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
